# The "Förändrad" (Changed) date in column C for every data row (2-210)
# is bumped forward by one day: 45177 -> 45178 (2023-09-08 -> 2023-09-09).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C210").Value = 45178
